# "making vector of profs"
# Adds two new working-hours log entries (rows 37 & 38) to Sheet1, recording
# two days of work on database manipulation / meshing professor data onto
# the main student dataframe.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 37: Monday, Jun 17 2024 -------------------------------------------------
$ws.Cells.Item(37, 1).Value = 45460            # Date
$ws.Cells.Item(37, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(37, 2).Value = "M"              # Day of week
$ws.Cells.Item(37, 2).NumberFormat = "d-mmm"
$ws.Cells.Item(37, 3).Value = 2                # Hours
$ws.Cells.Item(37, 5).Value = "database manipulation and "
$ws.Cells.Item(37, 5).WrapText = $true

# --- Row 38: Tuesday, Jun 18 2024 ------------------------------------------------
$ws.Cells.Item(38, 1).Value = 45461            # Date
$ws.Cells.Item(38, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(38, 2).Value = "T"              # Day of week
$ws.Cells.Item(38, 2).NumberFormat = "d-mmm"
$ws.Cells.Item(38, 3).Value = 2                # Hours

# Note: the "need to ..." aside (column G) was jotted down before the main
# note (column E) ended up being finalised, so it is written first here -
# this reproduces the original shared-string insertion order.
$ws.Cells.Item(38, 7).Value = "need to make length of student major minors equal to the length of unique student ids, then join. Probably can just use NA skip over profs who wont exist…"
$ws.Cells.Item(38, 5).Value = "trying to mesh profs onto main df_bsc table"
$ws.Cells.Item(38, 5).WrapText = $true

# Leave the selection on the newly-added last note, matching where the
# author's cursor ended up.
$ws.Range("E38").Select() | Out-Null
